$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2006
$ws.Range("I69").Value = 2006
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 6018
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -5144
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 2006
$ws.Range("I72").Value = 2006
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 18054
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -13686
$ws.Range("N72").ClearContents()
$ws.Range("H86").Value = 3322.5557
$ws.Range("I86").Value = 3337.4443
$ws.Range("J86").Value = 3307.6667
$ws.Range("K86").Value = 3337.4443
$ws.Range("L86").Value = 3307.6667
$ws.Range("M86").Value = -2214.4443
$ws.Range("N86").Value = -5553.6667
$ws.Range("H89").Value = 3322.5557
$ws.Range("I89").Value = 3337.4443
$ws.Range("J89").Value = 3307.6667
$ws.Range("K89").Value = 16687.2215
$ws.Range("L89").Value = 16538.3335
$ws.Range("M89").Value = -11071.2215
$ws.Range("N89").Value = -27770.3335
$ws.Range("H98").Value = 2333.2144
$ws.Range("I98").Value = 2055.4167
$ws.Range("K98").Value = 2055.4167
$ws.Range("M98").Value = -557.4167000000002
$ws.Range("H122").Value = 2333.2144
$ws.Range("I122").Value = 2055.4167
$ws.Range("K122").Value = 6166.250100000001
$ws.Range("M122").Value = -3716.250100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 396.33334
$ws.Range("I3").Value = 396.33334
$ws.Range("K3").Value = 396.33334
$ws.Range("M3").Value = -281.33334
$ws.Range("H32").Value = 23614.016
$ws.Range("J32").Value = 169436.42
$ws.Range("L32").Value = 169436.42
$ws.Range("N32").Value = -170010.42
$ws.Range("H61").Value = 7730.875
$ws.Range("I61").Value = 7721.143
$ws.Range("J61").Value = 7799
$ws.Range("K61").Value = 7721.143
$ws.Range("L61").Value = 7799
$ws.Range("M61").Value = -7509.143
$ws.Range("N61").Value = -8223
$ws.Range("H74").Value = 1901.9231
$ws.Range("I74").Value = 1673.125
$ws.Range("K74").Value = 1673.125
$ws.Range("M74").Value = -799.125
$ws.Range("H77").Value = 1901.9231
$ws.Range("I77").Value = 1673.125
$ws.Range("K77").Value = 8365.625
$ws.Range("M77").Value = -3997.625
$ws.Range("H94").Value = 30329.666
$ws.Range("J94").Value = 30329.666
$ws.Range("L94").Value = 30329.666
$ws.Range("N94").Value = -32131.666
$ws.Range("H132").Value = 1502.5938
$ws.Range("I132").Value = 1469.5333
$ws.Range("K132").Value = 4408.5999
$ws.Range("M132").Value = -1878.5999
$ws.Range("H136").Value = 7730.875
$ws.Range("I136").Value = 7721.143
$ws.Range("J136").Value = 7799
$ws.Range("K136").Value = 23163.429
$ws.Range("L136").Value = 23397
$ws.Range("M136").Value = -20613.429
$ws.Range("N136").Value = -28497

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("M8").Value = -860
$ws.Range("H86").Value = 2226.353
$ws.Range("I86").Value = 2064.25
$ws.Range("J86").Value = 2615.4
$ws.Range("K86").Value = 2064.25
$ws.Range("L86").Value = 2615.4
$ws.Range("M86").Value = -941.25
$ws.Range("N86").Value = -4861.4
$ws.Range("H89").Value = 2226.353
$ws.Range("I89").Value = 2064.25
$ws.Range("J89").Value = 2615.4
$ws.Range("K89").Value = 10321.25
$ws.Range("L89").Value = 13077
$ws.Range("M89").Value = -4705.25
$ws.Range("N89").Value = -24309

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2194.5
$ws.Range("I10").Value = 1926
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 1926
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -1787
$ws.Range("N10").Value = -3278
$ws.Range("H31").Value = 2831.8572
$ws.Range("I31").Value = 2011.5
$ws.Range("J31").Value = 3259.8696
$ws.Range("K31").Value = 2011.5
$ws.Range("L31").Value = 3259.8696
$ws.Range("M31").Value = -1716.5
$ws.Range("N31").Value = -3849.8696
$ws.Range("H34").Value = 2831.8572
$ws.Range("I34").Value = 2011.5
$ws.Range("J34").Value = 3259.8696
$ws.Range("K34").Value = 2011.5
$ws.Range("L34").Value = 3259.8696
$ws.Range("M34").Value = -1809.5
$ws.Range("N34").Value = -3663.8696
$ws.Range("H35").Value = 2508.4614
$ws.Range("I35").Value = 1015.7143
$ws.Range("J35").Value = 4250
$ws.Range("K35").Value = 1015.7143
$ws.Range("L35").Value = 4250
$ws.Range("M35").Value = -721.7143
$ws.Range("N35").Value = -4838
$ws.Range("H132").Value = 2296.5962
$ws.Range("I132").Value = 2215.9048
$ws.Range("J132").Value = 2635.5
$ws.Range("K132").Value = 6647.714399999999
$ws.Range("L132").Value = 7906.5
$ws.Range("M132").Value = -4117.714399999999
$ws.Range("N132").Value = -12966.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 330.57144
$ws.Range("I23").Value = 400
$ws.Range("J23").Value = 319
$ws.Range("K23").Value = 1200
$ws.Range("L23").Value = 957
$ws.Range("M23").Value = -965
$ws.Range("N23").Value = -1427
$ws.Range("H64").Value = 2816.2354
$ws.Range("J64").Value = 3069.077
$ws.Range("L64").Value = 9207.231
$ws.Range("N64").Value = -9747.231
$ws.Range("H67").Value = 2816.2354
$ws.Range("J67").Value = 3069.077
$ws.Range("L67").Value = 9207.231
$ws.Range("N67").Value = -11079.231
$ws.Range("H70").Value = 5616
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11685
$ws.Range("H73").Value = 5616
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -10908
$ws.Range("H117").Value = 5332.1177
$ws.Range("J117").Value = 5332.1177
$ws.Range("L117").Value = 15996.3531
$ws.Range("N117").Value = -22880.3531
$ws.Range("H125").Value = 9995
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 9995
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 29985
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -39825

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4147.522
$ws.Range("I80").Value = 3037.923
$ws.Range("J80").Value = 5590
$ws.Range("K80").Value = 3037.923
$ws.Range("L80").Value = 5590
$ws.Range("M80").Value = -2039.923
$ws.Range("N80").Value = -7586
$ws.Range("H83").Value = 4147.522
$ws.Range("I83").Value = 3037.923
$ws.Range("J83").Value = 5590
$ws.Range("K83").Value = 15189.615
$ws.Range("L83").Value = 27950
$ws.Range("M83").Value = -10197.615
$ws.Range("N83").Value = -37934
$ws.Range("H122").Value = 4608.636
$ws.Range("I122").Value = 4116.1665
$ws.Range("K122").Value = 12348.4995
$ws.Range("M122").Value = -9898.499500000002
$ws.Range("H132").Value = 3263.2334
$ws.Range("I132").Value = 3384.4583
$ws.Range("K132").Value = 10153.3749
$ws.Range("M132").Value = -7623.374899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2480.1538
$ws.Range("I16").Value = 2299.8
$ws.Range("J16").Value = 3081.3333
$ws.Range("K16").Value = 2299.8
$ws.Range("L16").Value = 3081.3333
$ws.Range("M16").Value = -2129.8
$ws.Range("N16").Value = -3421.3333
$ws.Range("H46").Value = 25843.166
$ws.Range("I46").Value = 40466.09
$ws.Range("K46").Value = 40466.09
$ws.Range("M46").Value = -40278.09
$ws.Range("H55").Value = 415.27777
$ws.Range("I55").Value = 218.41667
$ws.Range("J55").Value = 809
$ws.Range("K55").Value = 218.41667
$ws.Range("L55").Value = 809
$ws.Range("M55").Value = -45.41667000000001
$ws.Range("N55").Value = -1155
$ws.Range("H68").Value = 2261.5
$ws.Range("I68").Value = 1685.75
$ws.Range("J68").Value = 2837.25
$ws.Range("K68").Value = 1685.75
$ws.Range("L68").Value = 2837.25
$ws.Range("M68").Value = -936.75
$ws.Range("N68").Value = -4335.25
$ws.Range("H71").Value = 2261.5
$ws.Range("I71").Value = 1685.75
$ws.Range("J71").Value = 2837.25
$ws.Range("K71").Value = 8428.75
$ws.Range("L71").Value = 14186.25
$ws.Range("M71").Value = -4684.75
$ws.Range("N71").Value = -21674.25
$ws.Range("H122").Value = 13719.7
$ws.Range("I122").Value = 19449.834
$ws.Range("J122").Value = 5124.5
$ws.Range("K122").Value = 58349.50199999999
$ws.Range("L122").Value = 15373.5
$ws.Range("M122").Value = -55899.50199999999
$ws.Range("N122").Value = -20273.5
$ws.Range("H132").Value = 3271
$ws.Range("J132").Value = 4138.5625
$ws.Range("L132").Value = 12415.6875
$ws.Range("N132").Value = -17475.6875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3209.8125
$ws.Range("I81").Value = 3024.724
$ws.Range("J81").Value = 4999
$ws.Range("K81").Value = 6049.448
$ws.Range("L81").Value = 9998
$ws.Range("M81").Value = -4988.448
$ws.Range("N81").Value = -12120
$ws.Range("H84").Value = 3209.8125
$ws.Range("I84").Value = 3024.724
$ws.Range("J84").Value = 4999
$ws.Range("K84").Value = 30247.24
$ws.Range("L84").Value = 49990
$ws.Range("M84").Value = -24943.24
$ws.Range("N84").Value = -60598
$ws.Range("H107").Value = 26317068
$ws.Range("I107").Value = 1336.091
$ws.Range("K107").Value = 4008.273
$ws.Range("M107").Value = -2088.273
$ws.Range("H126").Value = 8492.277
$ws.Range("I126").Value = 2243.7273
$ws.Range("K126").Value = 6731.1819
$ws.Range("M126").Value = -4261.1819
$ws.Range("H132").Value = 2566871.8
$ws.Range("I132").Value = 1585674.8
$ws.Range("J132").Value = 5292419
$ws.Range("K132").Value = 4757024.4
$ws.Range("L132").Value = 15877257
$ws.Range("M132").Value = -4754494.4
$ws.Range("N132").Value = -15882317
